$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds numeric-looking strings (e.g. "589.10") that Excel would
# otherwise auto-convert into real numbers (losing the exact text / trailing zeros / precision).
# To keep them as literal text - matching the original inlineStr cells - briefly force a Text
# number format before assigning the value, then restore the default "Normal" style afterwards
# so no stray style index is left on the cell (matches original, unstyled data cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.057.19"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.130.32"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "589.10"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "147.86"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("D8").Value = "3.126.94"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +12.80%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +4.32%  "
$ws.Range("D14").Value = "37.56"
$ws.Range("E14").Value = "  +5.42%  "
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "3.648.95"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "63.880.63"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "3.128.12"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "467.83"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "14.38"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "13.28"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").Value = "82.44"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "9.01"
$ws.Range("E27").Value = "  +8.94%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "6.89"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "27.18"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").Value = "0.0₃0893"
$ws.Range("E34").Value = "  +10.42%  "
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "3.44"
$ws.Range("E37").Value = "  +12.27%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "50.96"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "8.71"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "0.0374"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "2.897.27"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").Value = "0.279"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "35.82"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "126.02"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "24.83"
$ws.Range("E51").Value = "  +0.66%  "

$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
